# Applies numeric data refresh to the Leve profit-tracking sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N), row by row,
# matching the upstream scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 108.818184
$ws.Range("I9").Value = 119.625
$ws.Range("K9").Value = 119.625
$ws.Range("M9").Value = 49.375

# Row 51
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3968

# Row 55
$ws.Range("H55").Value = 4898
$ws.Range("J55").Value = 4898
$ws.Range("L55").Value = 4898
$ws.Range("N55").Value = -5326

# Row 106
$ws.Range("H106").Value = 2799.5
$ws.Range("I106").Value = 1899
$ws.Range("J106").Value = 3700
$ws.Range("K106").Value = 1899
$ws.Range("L106").Value = 3700
$ws.Range("M106").Value = -1268
$ws.Range("N106").Value = -4962

# Row 123
$ws.Range("H123").Value = 49333
$ws.Range("J123").Value = 49333
$ws.Range("L123").Value = 49333
$ws.Range("N123").Value = -59133

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3999.8572
$ws.Range("I61").Value = 3999.8572
$ws.Range("K61").Value = 3999.8572
$ws.Range("M61").Value = -3787.8572

# Row 86
$ws.Range("H86").Value = 37500
$ws.Range("J86").Value = 37500
$ws.Range("L86").Value = 37500
$ws.Range("N86").Value = -39872

# Row 89
$ws.Range("H89").Value = 37500
$ws.Range("J89").Value = 37500
$ws.Range("L89").Value = 112500
$ws.Range("N89").Value = -124356

# Row 132
$ws.Range("H132").Value = 1250
$ws.Range("I132").Value = 1250
$ws.Range("K132").Value = 3750
$ws.Range("M132").Value = -1220

# Row 136
$ws.Range("H136").Value = 3999.8572
$ws.Range("I136").Value = 3999.8572
$ws.Range("K136").Value = 11999.5716
$ws.Range("M136").Value = -9449.571599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 635.4286
$ws.Range("I94").Value = 587.25
$ws.Range("K94").Value = 587.25
$ws.Range("M94").Value = -136.25

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2805.875
$ws.Range("I22").Value = 1899.5
$ws.Range("J22").Value = 3712.25
$ws.Range("K22").Value = 1899.5
$ws.Range("L22").Value = 3712.25
$ws.Range("M22").Value = -1549.5
$ws.Range("N22").Value = -4412.25

# Row 41
$ws.Range("H41").Value = 209870.23
$ws.Range("I41").Value = 9719.111000000001
$ws.Range("J41").Value = 299938.25
$ws.Range("K41").Value = 9719.111000000001
$ws.Range("L41").Value = 299938.25
$ws.Range("M41").Value = -9291.111000000001
$ws.Range("N41").Value = -300794.25

# Row 52
$ws.Range("H52").Value = 215831.67
$ws.Range("I52").Value = 44990
$ws.Range("J52").Value = 250000
$ws.Range("K52").Value = 44990
$ws.Range("L52").Value = 250000
$ws.Range("M52").Value = -44696
$ws.Range("N52").Value = -250588

# Row 95
$ws.Range("H95").Value = 19162.166
$ws.Range("J95").Value = 19162.166
$ws.Range("L95").Value = 19162.166
$ws.Range("N95").Value = -24654.166

# Row 132
$ws.Range("H132").Value = 4350
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

# Row 141
$ws.Range("H141").Value = 80983.5
$ws.Range("J141").Value = 80983.5
$ws.Range("L141").Value = 80983.5
$ws.Range("N141").Value = -91343.5

$ws = $wb.Worksheets.Item("CUL")
# Row 28
$ws.Range("H28").Value = 1700
$ws.Range("I28").Value = 1700
$ws.Range("K28").Value = 5100
$ws.Range("M28").Value = -4868

# Row 139
$ws.Range("H139").Value = 3867.6667
$ws.Range("I139").Value = 2014
$ws.Range("K139").Value = 6042
$ws.Range("M139").Value = -902

# Row 140
$ws.Range("H140").Value = 4998
$ws.Range("I140").Value = 4997
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 14991
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -9811
$ws.Range("N140").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 36.11111
$ws.Range("J2").Value = 60.2
$ws.Range("L2").Value = 60.2
$ws.Range("N2").Value = -286.2

# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# Row 36
$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Row 80
$ws.Range("H80").Value = 3529.4443
$ws.Range("I80").Value = 3279.8
$ws.Range("J80").Value = 3841.5
$ws.Range("K80").Value = 3279.8
$ws.Range("L80").Value = 3841.5
$ws.Range("M80").Value = -2281.8
$ws.Range("N80").Value = -5837.5

# Row 83
$ws.Range("H83").Value = 3529.4443
$ws.Range("I83").Value = 3279.8
$ws.Range("J83").Value = 3841.5
$ws.Range("K83").Value = 16399
$ws.Range("L83").Value = 19207.5
$ws.Range("M83").Value = -11407
$ws.Range("N83").Value = -29191.5

# Row 102
$ws.Range("H102").Value = 1944.3
$ws.Range("I102").Value = 1944.3
$ws.Range("K102").Value = 1944.3
$ws.Range("M102").Value = -322.3

# Row 122
$ws.Range("H122").Value = 2773.0557
$ws.Range("I122").Value = 2120.4375
$ws.Range("K122").Value = 6361.3125
$ws.Range("M122").Value = -3911.3125

# Row 126
$ws.Range("H126").Value = 3453.7646
$ws.Range("I126").Value = 2594.6
$ws.Range("K126").Value = 7783.799999999999
$ws.Range("M126").Value = -5313.799999999999

# Row 129
$ws.Range("H129").Value = 70750
$ws.Range("J129").Value = 70750
$ws.Range("L129").Value = 70750
$ws.Range("N129").Value = -80750

# Row 132
$ws.Range("H132").Value = 4967.25
$ws.Range("I132").Value = 5221.1665
$ws.Range("K132").Value = 15663.4995
$ws.Range("M132").Value = -13133.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8796.799999999999
$ws.Range("I7").Value = 7997.6665
$ws.Range("K7").Value = 7997.6665
$ws.Range("M7").Value = -7885.6665

# Row 40
$ws.Range("H40").Value = 4050.4285
$ws.Range("I40").Value = 4050.4285
$ws.Range("K40").Value = 4050.4285
$ws.Range("M40").Value = -3914.4285

# Row 88
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 20000
$ws.Range("K88").Value = 20000
$ws.Range("M88").Value = -19572

# Row 91
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 20000
$ws.Range("K91").Value = 20000
$ws.Range("M91").Value = -18518

# Row 95
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

# Row 126
$ws.Range("H126").Value = 8796.799999999999
$ws.Range("I126").Value = 7997.6665
$ws.Range("K126").Value = 23992.9995
$ws.Range("M126").Value = -21522.9995

# Row 136
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("WVR")
# Row 48
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -51137

# Row 98
$ws.Range("H98").Value = 43737.5
$ws.Range("J98").Value = 43737.5
$ws.Range("L98").Value = 43737.5
$ws.Range("N98").Value = -49727.5

# Row 127
$ws.Range("H127").Value = 125000
$ws.Range("J127").Value = 125000
$ws.Range("L127").Value = 125000
$ws.Range("N127").Value = -134920

# Row 136
$ws.Range("H136").Value = 1519.8
$ws.Range("I136").Value = 849.5
$ws.Range("J136").Value = 1966.6666
$ws.Range("K136").Value = 2548.5
$ws.Range("L136").Value = 5899.9998
$ws.Range("M136").Value = 1.5
$ws.Range("N136").Value = -10999.9998
